$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (date) values between row 3 and row 5
$ws.Range("D3").Value = 44257
$ws.Range("D5").Value = 44250

# Swap the Volumen values between row 3 and row 5
$ws.Range("M3").Value = 100
$ws.Range("M5").Value = 200
